$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format so the numeric-looking strings keep their
# exact textual representation (trailing zeros, precision) instead of being
# coerced into floating point numbers when the .Value is assigned.
$ranges = @("D2:D12", "D14:D24", "D41:D51")
foreach ($r in $ranges) {
    $ws.Range($r).NumberFormat = "@"
}

# Updated coin prices (column D) scraped on 2022-12-25
$ws.Range("D2").Value = "241.99"
$ws.Range("D3").Value = "22.90"
$ws.Range("D4").Value = "5.386"
$ws.Range("D5").Value = "0.05941"
$ws.Range("D6").Value = "3.398"
$ws.Range("D7").Value = "6.452"
$ws.Range("D8").Value = "0.8046"
$ws.Range("D9").Value = "0.9134"
$ws.Range("D10").Value = "0.1410"
$ws.Range("D11").Value = "0.07416"
$ws.Range("D12").Value = "0.03281"
$ws.Range("D14").Value = "0.09315"
$ws.Range("D15").Value = "3.863"
$ws.Range("D16").Value = "0.001572"
$ws.Range("D17").Value = "0.04525"
$ws.Range("D18").Value = "0.0005946"
$ws.Range("D19").Value = "0.006086"
$ws.Range("D20").Value = "0.004992"
$ws.Range("D21").Value = "0.007494"
$ws.Range("D22").Value = "0.0009865"
$ws.Range("D23").Value = "0.00007808"
$ws.Range("D24").Value = "3.611"
$ws.Range("D41").Value = "0.006229"
$ws.Range("D42").Value = "0.1062"
$ws.Range("D43").Value = "0.002803"
$ws.Range("D44").Value = "0.007214"
$ws.Range("D45").Value = "0.00005186"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").Value = "0.0005805"
$ws.Range("D48").Value = "0.9589"
$ws.Range("D49").Value = "0.002263"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").Value = "0.0002002"

# Restore the default (Normal) style so no visible formatting/number-format
# change is left behind on the cells -- only their text content changes.
foreach ($r in $ranges) {
    $ws.Range($r).Style = "Normal"
}
